# Commit: Thu, Jul 09, 2020  8:05:04 PM
#
# This edit does two things to the deck:
#
# 1) Three tables (on slides 14, 15 and 16) get their table style switched
#    from the deck's custom "Table_0" style ({4011D032-51E3-4C0D-8E9D-
#    A4E59DF9DBC6}) to the PowerPoint built-in "No Style, Table Grid"
#    style ({9F544921-FA93-46B4-9268-91FF1CF3EC88}).
#
# 2) The presentation's theme colour scheme (the one driving the slide
#    master / all slides) is switched from the "Integral" / "Red Violet"
#    palette to the standard "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1) Re-style the three tables -----------------------------------------

$tableStyleId = "{9F544921-FA93-46B4-9268-91FF1CF3EC88}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($tableStyleId, $true)
        }
    }
}

# --- 2) Swap the active theme colours to the standard Office palette ------

$officeThemeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColors[$i - 1]
}
